# Apply the recorded change to the Denmark Superliga 2023-2024 sheet:
#  1) Rows 63 and 64 had their match data (columns F..V) swapped back to the
#     correct order (Nordsjaelland-Odense / Vejle-Hvidovre IF).
#  2) A new row 98 was appended with the Randers FC vs Vejle match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the F:V contents of rows 63 and 64 -----------------------

$row63 = $ws.Range("F63:V63").Value()
$row64 = $ws.Range("F64:V64").Value()

$ws.Range("F63:V63").Value = $row64
$ws.Range("F64:V64").Value = $row63

# --- 2) Append new row 98 ----------------------------------------------

$r = 98

$ws.Cells.Item($r, 1).Value = 97
$ws.Cells.Item($r, 2).Value = "denmark"
$ws.Cells.Item($r, 3).Value = "superliga"
$ws.Cells.Item($r, 4).Value = "2023-2024"
$ws.Cells.Item($r, 5).Value = 45261.79166666666
$ws.Cells.Item($r, 6).Value = "Randers FC"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = "Vejle"
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 1.98
$ws.Cells.Item($r, 11).Value = "26/11/2023 16:13"
$ws.Cells.Item($r, 12).Value = 2.04
$ws.Cells.Item($r, 13).Value = "01/12/2023 18:52"
$ws.Cells.Item($r, 14).Value = 3.55
$ws.Cells.Item($r, 15).Value = "26/11/2023 16:13"
$ws.Cells.Item($r, 16).Value = 3.47
$ws.Cells.Item($r, 17).Value = "01/12/2023 18:52"
$ws.Cells.Item($r, 18).Value = 3.88
$ws.Cells.Item($r, 19).Value = "26/11/2023 16:13"
$ws.Cells.Item($r, 20).Value = 3.89
$ws.Cells.Item($r, 21).Value = "01/12/2023 18:52"
$ws.Cells.Item($r, 22).Value = "https://www.betexplorer.com/football/denmark/superliga/randers-fc-vejle/0zXFLXVE/"

# Match the "Indice" (A) column styling (bold, centered/top aligned, thin
# border) and the "data_partida" (E) column's datetime number format, as
# used by every other data row in the sheet.
$ws.Cells.Item($r, 1).Font.Bold = $true
$ws.Cells.Item($r, 1).HorizontalAlignment = -4108
$ws.Cells.Item($r, 1).VerticalAlignment = -4160
$ws.Cells.Item($r, 1).Borders.LineStyle = 1

$ws.Cells.Item($r, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
